# Update "docs/epexspot_prices.xlsx"
#  - Sheet "Prix Spot": add a new day column BD (08-aug) with its 24 hourly prices
#  - Sheet "Gaz": append a new row (2025-08-06 / 32.175)
#  - Sheet "CO2": append a new row (2025-08-06 / 70.23999999999999)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Prix Spot" - new column BD ("08-aug")
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the previous day column (BC) onto the new one (BD)
# so the header keeps the same bold/border/centered style.
$wsPrix.Range("BC1").Copy()
$wsPrix.Range("BD1").PasteSpecial(-4122)  # xlPasteFormats

$wsPrix.Range("BD1").Value = "08-aug"

$prixValues = @{
    2  = 72.13
    3  = 64.44
    4  = 53.67
    5  = 47.28
    6  = 52.29
    7  = 70.28
    8  = 73
    9  = 84.5
    10 = 78.69
    11 = 68.65000000000001
    12 = 41.13
    13 = 8.460000000000001
    14 = 0
    15 = -0.01
    16 = -0.02
    17 = -0.01
    18 = 10.49
    19 = 49.77
    20 = 68.62
    21 = 104.07
    22 = 96.16
    23 = 77.17
    24 = 85.7
    25 = 86.2
}

foreach ($row in $prixValues.Keys) {
    $wsPrix.Cells.Item($row, 56).Value = $prixValues[$row]  # column 56 = BD
}

# ---------------------------------------------------------------------------
# Sheet 2: "Gaz" - append row 53
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force the date to be stored as plain text (matching the other rows) rather
# than letting Excel auto-convert "2025-08-06" into a date serial number.
$wsGaz.Range("A53").NumberFormat = "@"
$wsGaz.Range("A53").Value = "2025-08-06"
$wsGaz.Range("A53").ClearFormats()

$wsGaz.Range("B53").Value = 32.175

# ---------------------------------------------------------------------------
# Sheet 3: "CO2" - append row 53
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A53").NumberFormat = "@"
$wsCo2.Range("A53").Value = "2025-08-06"
$wsCo2.Range("A53").ClearFormats()

$wsCo2.Range("B53").Value = 70.23999999999999
